$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 209; this shifts the existing rows 209-244
# down to 210-245, preserving all of their data and formatting.
$ws.Rows.Item(209).Insert()

# Populate the newly inserted row 209 with the new weekly record.
$ws.Cells.Item(209, 1).Value = 8
$ws.Cells.Item(209, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(209, 3).Value = "Coquimbo"
$ws.Cells.Item(209, 4).Value = 44476
$ws.Cells.Item(209, 5).Value = 4
$ws.Cells.Item(209, 6).Value = 100114001
$ws.Cells.Item(209, 7).Value = "Papa"
$ws.Cells.Item(209, 8).Value = "Cardinal"
$ws.Cells.Item(209, 9).Value = "1a (cosecha)"
$ws.Cells.Item(209, 10).Value = 2000
$ws.Cells.Item(209, 11).Value = 13000
$ws.Cells.Item(209, 12).Value = 14000
$ws.Cells.Item(209, 13).Value = 13500
$ws.Cells.Item(209, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(209, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(209, 16).Value = 540
$ws.Cells.Item(209, 17).Value = 25
$ws.Cells.Item(209, 18).Value = "Hortaliza"

# Match the numeric date formatting style used by the other date cells
# in column D (style index 2 in the original workbook).
$ws.Cells.Item(209, 4).NumberFormat = $ws.Cells.Item(210, 4).NumberFormat
